# Actualización automática 2025-06-16 13:01:14
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new column before column F ("GRIFERIAS") for the new "GRANITO" group.
# This shifts the existing F:N columns one place to the right (to G:O).
$ws.Columns("F:F").Insert()

# Match the formatting of the neighboring "FREGADEROS DE COCINA" column (E) for
# the header/value/summary styles (s=1 / s=2 / s=4).
$ws.Range("E1:E7").Copy()
$ws.Range("F1:F7").PasteSpecial(-4122)

# New header for the inserted column.
$ws.Cells.Item(1, 6).Value = "GRANITO"

# Fill the inserted column's data rows (2-6) with 0, matching the other group columns.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Row 7 summary cell for the new column.
$ws.Cells.Item(7, 6).Value = "0 de 5"

# Append three brand-new columns after the existing last column (now O = 15),
# matching the formatting of column O ("SAL SOLUBLE").
$ws.Range("O1:O7").Copy()
$ws.Range("P1:R7").PasteSpecial(-4122)

$headers = @("NO RESURTIBLES", "PANELES PVC", "PANELES PU")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = 16 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    for ($r = 2; $r -le 6; $r++) {
        $ws.Cells.Item($r, $col).Value = 0
    }
    $ws.Cells.Item(7, $col).Value = "0 de 5"
}

# Column widths per the target layout (columns F..R). The ColumnWidth setter
# round-trips through a pixel conversion that adds ~0.8333 to the stored
# <col width>, so compensate by subtracting it to land on the exact value.
$padding = 0.8333333333333334
$widths = @{
    6  = 13
    7  = 15
    8  = 14
    9  = 13
    10 = 9
    11 = 25
    12 = 24
    13 = 17
    14 = 26
    15 = 17
    16 = 20
    17 = 17
    18 = 16
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - $padding
}
